$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted numeric-looking values (e.g. "1.003") stay as text, not numbers
$ws.Range('D2').Value = '30.477.07'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '1.911.82'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.40'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4786'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2839'
$ws.Range('E8').Value = '  -3.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06698'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.85'
$ws.Range('E10').Value = '  -2.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '101.88'
$ws.Range('E11').Value = '  -3.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07711'
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').Value = '1.917.65'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.204'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6699'
$ws.Range('E15').Value = '  -4.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '267.10'
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('D17').Value = '30.505.99'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007467'
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.67'
$ws.Range('E20').Value = '  -3.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.418'
$ws.Range('E21').Value = '  -1.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.294'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.379'
$ws.Range('E24').Value = '  -3.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '166.86'
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.19'
$ws.Range('E26').Value = '  -2.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.063'
$ws.Range('E27').Value = '  -4.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.391'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  -4.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.620'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.515'
$ws.Range('E31').Value = '  -2.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.214'
$ws.Range('E32').Value = '  -3.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04722'
$ws.Range('E33').Value = '  -2.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7259'
$ws.Range('E34').Value = '  -3.78%  '
$ws.Range('E35').Value = '  -4.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.722'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01913'
$ws.Range('E37').Value = '  -4.13%  '
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.289'
$ws.Range('E39').Value = '  -3.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '74.82'
$ws.Range('E40').Value = '  -3.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.969'
$ws.Range('E41').Value = '  -6.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8601'
$ws.Range('E42').Value = '  -4.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '105.05'
$ws.Range('E43').Value = '  -2.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4261'
$ws.Range('E44').Value = '  -3.18%  '
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.397'
$ws.Range('E46').Value = '  -4.65%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '923.61'
$ws.Range('E47').Value = '  -8.06%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1201'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.77'
$ws.Range('E49').Value = '  -3.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.817'
$ws.Range('E50').Value = '  -5.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05759'
$ws.Range('E51').Value = '  +0.52%  '
